$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 28.27788229977268
$ws.Range("D2").Value = -1.092117700227323
$ws.Range("E2").Value = 1.192721071149818
$ws.Range("C3").Value = 29.19175924173338
$ws.Range("D3").Value = -0.3482407582666198
$ws.Range("E3").Value = 0.1212716257181103
$ws.Range("C4").Value = 30.09155300362499
$ws.Range("D4").Value = 0.541553003624994
$ws.Range("E4").Value = 0.2932796557352528
$ws.Range("C5").Value = 29.78267262292169
$ws.Range("D5").Value = 0.03267262292168738
$ws.Range("E5").Value = 0.001067500288582771
$ws.Range("C6").Value = 29.75900732942098
$ws.Range("D6").Value = -0.08099267057901827
$ws.Range("E6").Value = 0.006559812687521371
$ws.Range("C7").Value = 29.84561044431663
$ws.Range("D7").Value = 0.03561044431662808
$ws.Range("E7").Value = 0.001268103744427669
$ws.Range("C8").Value = 29.62121671700346
$ws.Range("D8").Value = -0.2987832829965384
$ws.Range("E8").Value = 0.08927145019818958
$ws.Range("C9").Value = 29.83432592272917
$ws.Range("D9").Value = -0.1456740772708294
$ws.Range("E9").Value = 0.02122093678870757
$ws.Range("C10").Value = 29.88982670961775
$ws.Range("D10").Value = -0.1501732903822486
$ws.Range("E10").Value = 0.02255201714423116
$ws.Range("C11").Value = 29.78883031723733
$ws.Range("D11").Value = -0.4211696827626668
$ws.Range("E11").Value = 0.1773839016784053
$ws.Range("C12").Value = 30.16589280729863
$ws.Range("D12").Value = -0.05410719270137321
$ws.Range("E12").Value = 0.002927588302023534
$ws.Range("C13").Value = 29.82036427603519
$ws.Range("D13").Value = -0.5596357239648135
$ws.Range("E13").Value = 0.313192143537621
$ws.Range("C14").Value = 30.05698656662453
$ws.Range("D14").Value = -0.3830134333754707
$ws.Range("E14").Value = 0.1466992901460661
$ws.Range("C15").Value = 30.19023929771748
$ws.Range("D15").Value = -0.2897607022825248
$ws.Range("E15").Value = 0.08396126458726198
$ws.Range("C16").Value = 29.96611093429973
$ws.Range("D16").Value = -0.7238890657002663
$ws.Range("E16").Value = 0.5240153794404045
$ws.Range("C17").Value = 29.95353145807043
$ws.Range("D17").Value = -0.7964685419295705
$ws.Range("E17").Value = 0.634362138283416
$ws.Range("C18").Value = 30.54948913695124
$ws.Range("D18").Value = -0.3905108630487568
$ws.Range("E18").Value = 0.1524987341590849
$ws.Range("C19").Value = 30.77586982982714
$ws.Range("D19").Value = -0.1741301701728553
$ws.Range("E19").Value = 0.03032131616442753
$ws.Range("C20").Value = 30.70160810776703
$ws.Range("D20").Value = -0.3183918922329703
$ws.Range("E20").Value = 0.1013733970396914
$ws.Range("C21").Value = 30.78623034825427
$ws.Range("D21").Value = -0.3337696517457296
$ws.Range("E21").Value = 0.1114021804264656
$ws.Range("C22").Value = 31.27141867503842
$ws.Range("D22").Value = -0.008581324961582482
$ws.Range("E22").Value = 0.00007363913809627858
$ws.Range("C23").Value = 31.3953411648987
$ws.Range("D23").Value = 0.01534116489870385
$ws.Range("E23").Value = 0.0002353513404492232
$ws.Range("C24").Value = 31.58726286396578
$ws.Range("D24").Value = 0.00726286396577791
$ws.Range("E24").Value = 0.00005274919298539523
$ws.Range("C25").Value = 31.5674891895968
$ws.Range("D25").Value = -0.08251081040319974
$ws.Range("E25").Value = 0.006808033833392775
$ws.Range("C26").Value = 31.90485810887399
$ws.Range("D26").Value = 0.02485810887398898
$ws.Range("E26").Value = 0.0006179255767910898
$ws.Range("C27").Value = 32.29400698910788
$ws.Range("D27").Value = 0.01400698910787668
$ws.Range("E27").Value = 0.0001961957438681759
$ws.Range("C28").Value = 32.24138097474844
$ws.Range("D28").Value = -0.2086190252515649
$ws.Range("E28").Value = 0.04352189769691309
$ws.Range("C29").Value = 33.52817252141637
$ws.Range("D29").Value = 0.6781725214163714
$ws.Range("E29").Value = 0.4599179688042387
$ws.Range("C30").Value = 33.28296385789197
$ws.Range("D30").Value = 0.3829638578919727
$ws.Range("E30").Value = 0.1466613164515031
$ws.Range("C31").Value = 33.3784015109037
$ws.Range("D31").Value = 0.2784015109036986
$ws.Range("E31").Value = 0.0775074012734622
$ws.Range("C32").Value = 33.25763537194656
$ws.Range("D32").Value = -0.1423646280534356
$ws.Range("E32").Value = 0.02026768732079307
$ws.Range("C33").Value = 33.94127858902861
$ws.Range("D33").Value = 0.2412785890286031
$ws.Range("E33").Value = 0.05821535752363356
$ws.Range("C34").Value = 34.80115297861217
$ws.Range("D34").Value = 0.7011529786121642
$ws.Range("E34").Value = 0.49161549941671
$ws.Range("C35").Value = 34.71223657580823
$ws.Range("D35").Value = 0.3122365758082282
$ws.Range("E35").Value = 0.09749167927244741
$ws.Range("C36").Value = 35.20854410149703
$ws.Range("D36").Value = 0.3085441014970343
$ws.Range("E36").Value = 0.09519946256861218
$ws.Range("C37").Value = 35.00070622077768
$ws.Range("D37").Value = -0.2992937792223174
$ws.Range("E37").Value = 0.08957676628117725
$ws.Range("C38").Value = 35.39205346563158
$ws.Range("D38").Value = -0.3079465343684191
$ws.Range("E38").Value = 0.09483106802951995
$ws.Range("C39").Value = 35.98064298204398
$ws.Range("D39").Value = -0.3193570179560155
$ws.Range("E39").Value = 0.1019889049177588
$ws.Range("C40").Value = 36.42452701316659
$ws.Range("D40").Value = -0.3754729868334081
$ws.Range("E40").Value = 0.1409799638416007
$ws.Range("C41").Value = 37.51103321642204
$ws.Range("D41").Value = 0.2110332164220452
$ws.Range("E41").Value = 0.04453501843343377
$ws.Range("C42").Value = 38.35515710999172
$ws.Range("D42").Value = 0.455157109991724
$ws.Range("E42").Value = 0.2071679947760183
$ws.Range("C43").Value = 38.90910449529751
$ws.Range("D43").Value = 0.4091044952975125
$ws.Range("E43").Value = 0.1673664880726324
$ws.Range("C44").Value = 39.76024284485575
$ws.Range("D44").Value = 0.8602428448557546
$ws.Range("E44").Value = 0.7400177521255219
$ws.Range("C45").Value = 40.40243169098477
$ws.Range("D45").Value = 1.002431690984771
$ws.Range("E45").Value = 1.004869295090588
$ws.Range("C46").Value = 40.64866100100531
$ws.Range("D46").Value = 0.7486610010053099
$ws.Range("E46").Value = 0.5604932944262726
$ws.Range("C47").Value = 39.58036945558321
$ws.Range("D47").Value = -0.5196305444167919
$ws.Range("E47").Value = 0.2700159026908915
$ws.Range("C48").Value = 40.98050448704024
$ws.Range("D48").Value = 0.3805044870402341
$ws.Range("E48").Value = 0.1447836646577517
$ws.Range("C49").Value = 41.08319000266954
$ws.Range("D49").Value = 0.1831900026695408
$ws.Range("E49").Value = 0.03355857707806636
$ws.Range("C50").Value = 40.76523899808499
$ws.Range("D50").Value = -0.4347610019150139
$ws.Range("E50").Value = 0.1890171287861467
$ws.Range("C51").Value = 40.21752257365866
$ws.Range("D51").Value = -1.282477426341345
$ws.Range("E51").Value = 1.644748349075119
$ws.Range("C52").Value = -2.717463598228047
$ws.Range("E52").Value = 10.9596818406561
$ws.Range("E53").Value = 0.2191936368131221
